# Slide 3, shape "Скругленная прямоугольная выноска 10" ("Сохранить пример"):
# put the cursor at the end of the existing text and press Enter three
# times, leaving three new empty (centered) paragraphs below the label,
# each inheriting the same red / 12pt formatting as the original text.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item("Скругленная прямоугольная выноска 10")

$tr = $shp.TextFrame.TextRange
$tr.InsertAfter("`r`r`r")
